{"js": "// Replace each math-expression cell's text with its new value, in document order.\n// Each \"before\" string is unique in the document, so a literal (non-wildcard) search\n// locates exactly the one cell to update.\nconst replacements = [\n  [\"93-76=\", \"91-9=\"],\n  [\"92-28=\", \"70+12=\"],\n  [\"45+15=\", \"7+47=\"],\n  [\"46-22=\", \"87-26=\"],\n  [\"7+69=\", \"42+8=\"],\n  [\"58-52=\", \"49-14=\"],\n  [\"31+42=\", \"62-5=\"],\n  [\"71+8=\", \"83-81=\"],\n  [\"1+92=\", \"70-6=\"],\n  [\"98-93=\", \"85-37=\"],\n  [\"6+63=\", \"97-66=\"],\n  [\"48+38=\", \"94-35=\"],\n  [\"32-28=\", \"85-21=\"],\n  [\"76-72=\", \"16+43=\"],\n  [\"90-12=\", \"91-19=\"],\n  [\"29+15=\", \"73-52=\"],\n  [\"25+36=\", \"54+24=\"],\n  [\"88-34=\", \"65-38=\"],\n  [\"51+7=\", \"23+30=\"],\n  [\"56+6=\", \"90-45=\"],\n  [\"39+19=\", \"72-69=\"],\n  [\"41+51=\", \"13+29=\"],\n  [\"75-19=\", \"62-12=\"],\n  [\"42+41=\", \"81-79=\"],\n  [\"40+5=\", \"63+0=\"],\n  [\"86-58=\", \"99-29=\"],\n  [\"64+8=\", \"85-22=\"],\n  [\"80-47=\", \"72+9=\"],\n  [\"2+62=\", \"32+9=\"],\n  [\"85-34=\", \"36+57=\"],\n  [\"29+3=\", \"35+49=\"],\n  [\"56+26=\", \"23+36=\"],\n  [\"67-18=\", \"60-28=\"],\n  [\"29+0=\", \"30-2=\"],\n  [\"53+21=\", \"7+6=\"],\n  [\"78-55=\", \"84-81=\"],\n  [\"86-67=\", \"38+10=\"],\n  [\"88-15=\", \"41-33=\"],\n  [\"65+29=\", \"61+11=\"],\n  [\"59+40=\", \"97-87=\"],\n  [\"66-6=\", \"57-39=\"],\n  [\"36+1=\", \"59-36=\"],\n  [\"52-0=\", \"1+11=\"],\n  [\"54+40=\", \"96-10=\"],\n  [\"19+63=\", \"73-59=\"],\n  [\"44+17=\", \"0+63=\"],\n  [\"12+1=\", \"7+65=\"],\n  [\"30+57=\", \"64+16=\"],\n  [\"5+6=\", \"62+34=\"],\n  [\"34+31=\", \"7+77=\"],\n  [\"40+20=\", \"56+42=\"],\n  [\"59-12=\", \"21-16=\"],\n  [\"44-41=\", \"73+10=\"],\n  [\"17+22=\", \"27-18=\"],\n  [\"68-2=\", \"91-39=\"],\n  [\"8+19=\", \"73-68=\"],\n  [\"59+1=\", \"15-10=\"],\n  [\"94+0=\", \"81-28=\"],\n  [\"6+3=\", \"0+67=\"],\n  [\"23-13=\", \"35+53=\"],\n  [\"76+20=\", \"42+17=\"],\n  [\"25+5=\", \"87-65=\"],\n  [\"84-39=\", \"67-37=\"],\n  [\"28+48=\", \"24+64=\"],\n  [\"27+3=\", \"63+5=\"],\n  [\"58+11=\", \"55-49=\"],\n  [\"65-24=\", \"51+26=\"],\n  [\"61+29=\", \"42+50=\"],\n  [\"77-75=\", \"45-1=\"],\n  [\"42+22=\", \"25+46=\"],\n  [\"52-17=\", \"57-32=\"],\n  [\"33-12=\", \"17-9=\"],\n  [\"96-26=\", \"83-48=\"],\n  [\"23+34=\", \"31+22=\"],\n  [\"27+24=\", \"84-82=\"],\n  [\"37+0=\", \"50-21=\"],\n  [\"75+2=\", \"96-25=\"],\n  [\"52-10=\", \"81+4=\"],\n  [\"42+9=\", \"9+12=\"],\n  [\"89-52=\", \"76-36=\"],\n  [\"20+1=\", \"23-19=\"],\n  [\"67-15=\", \"17+14=\"],\n  [\"71-44=\", \"0+70=\"],\n  [\"82-75=\", \"69-56=\"],\n  [\"2+10=\", \"1+90=\"],\n  [\"51+41=\", \"83-44=\"],\n  [\"16+45=\", \"90-53=\"],\n  [\"57-13=\", \"94-16=\"],\n  [\"56-41=\", \"93-49=\"],\n  [\"59+35=\", \"4+62=\"],\n  [\"74-28=\", \"6+19=\"],\n  [\"66-59=\", \"6+80=\"],\n  [\"13+26=\", \"20+39=\"],\n  [\"94-12=\", \"2+86=\"],\n  [\"78-23=\", \"78-66=\"],\n  [\"41-18=\", \"73-19=\"],\n  [\"72+8=\", \"40+36=\"],\n  [\"78+3=\", \"81-63=\"],\n  [\"68-9=\", \"12+52=\"],\n  [\"52+0=\", \"40-39=\"],\n];\n\nfor (const [before, after] of replacements) {\n  const results = context.document.body.search(before, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(`Text not found: ${before}`);\n  }\n\n  for (const item of results.items) {\n    item.insertText(after, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Replace each math-expression cell's text with its new value, in document order.\n# Each \"before\" string is unique in the document, so Find/Replace locates exactly\n# the one cell to update; no formatting is touched since only the text changes.\n$d = $word.ActiveDocument\n\n$pairs = @(\n    @(\"93-76=\", \"91-9=\"),\n    @(\"92-28=\", \"70+12=\"),\n    @(\"45+15=\", \"7+47=\"),\n    @(\"46-22=\", \"87-26=\"),\n    @(\"7+69=\", \"42+8=\"),\n    @(\"58-52=\", \"49-14=\"),\n    @(\"31+42=\", \"62-5=\"),\n    @(\"71+8=\", \"83-81=\"),\n    @(\"1+92=\", \"70-6=\"),\n    @(\"98-93=\", \"85-37=\"),\n    @(\"6+63=\", \"97-66=\"),\n    @(\"48+38=\", \"94-35=\"),\n    @(\"32-28=\", \"85-21=\"),\n    @(\"76-72=\", \"16+43=\"),\n    @(\"90-12=\", \"91-19=\"),\n    @(\"29+15=\", \"73-52=\"),\n    @(\"25+36=\", \"54+24=\"),\n    @(\"88-34=\", \"65-38=\"),\n    @(\"51+7=\", \"23+30=\"),\n    @(\"56+6=\", \"90-45=\"),\n    @(\"39+19=\", \"72-69=\"),\n    @(\"41+51=\", \"13+29=\"),\n    @(\"75-19=\", \"62-12=\"),\n    @(\"42+41=\", \"81-79=\"),\n    @(\"40+5=\", \"63+0=\"),\n    @(\"86-58=\", \"99-29=\"),\n    @(\"64+8=\", \"85-22=\"),\n    @(\"80-47=\", \"72+9=\"),\n    @(\"2+62=\", \"32+9=\"),\n    @(\"85-34=\", \"36+57=\"),\n    @(\"29+3=\", \"35+49=\"),\n    @(\"56+26=\", \"23+36=\"),\n    @(\"67-18=\", \"60-28=\"),\n    @(\"29+0=\", \"30-2=\"),\n    @(\"53+21=\", \"7+6=\"),\n    @(\"78-55=\", \"84-81=\"),\n    @(\"86-67=\", \"38+10=\"),\n    @(\"88-15=\", \"41-33=\"),\n    @(\"65+29=\", \"61+11=\"),\n    @(\"59+40=\", \"97-87=\"),\n    @(\"66-6=\", \"57-39=\"),\n    @(\"36+1=\", \"59-36=\"),\n    @(\"52-0=\", \"1+11=\"),\n    @(\"54+40=\", \"96-10=\"),\n    @(\"19+63=\", \"73-59=\"),\n    @(\"44+17=\", \"0+63=\"),\n    @(\"12+1=\", \"7+65=\"),\n    @(\"30+57=\", \"64+16=\"),\n    @(\"5+6=\", \"62+34=\"),\n    @(\"34+31=\", \"7+77=\"),\n    @(\"40+20=\", \"56+42=\"),\n    @(\"59-12=\", \"21-16=\"),\n    @(\"44-41=\", \"73+10=\"),\n    @(\"17+22=\", \"27-18=\"),\n    @(\"68-2=\", \"91-39=\"),\n    @(\"8+19=\", \"73-68=\"),\n    @(\"59+1=\", \"15-10=\"),\n    @(\"94+0=\", \"81-28=\"),\n    @(\"6+3=\", \"0+67=\"),\n    @(\"23-13=\", \"35+53=\"),\n    @(\"76+20=\", \"42+17=\"),\n    @(\"25+5=\", \"87-65=\"),\n    @(\"84-39=\", \"67-37=\"),\n    @(\"28+48=\", \"24+64=\"),\n    @(\"27+3=\", \"63+5=\"),\n    @(\"58+11=\", \"55-49=\"),\n    @(\"65-24=\", \"51+26=\"),\n    @(\"61+29=\", \"42+50=\"),\n    @(\"77-75=\", \"45-1=\"),\n    @(\"42+22=\", \"25+46=\"),\n    @(\"52-17=\", \"57-32=\"),\n    @(\"33-12=\", \"17-9=\"),\n    @(\"96-26=\", \"83-48=\"),\n    @(\"23+34=\", \"31+22=\"),\n    @(\"27+24=\", \"84-82=\"),\n    @(\"37+0=\", \"50-21=\"),\n    @(\"75+2=\", \"96-25=\"),\n    @(\"52-10=\", \"81+4=\"),\n    @(\"42+9=\", \"9+12=\"),\n    @(\"89-52=\", \"76-36=\"),\n    @(\"20+1=\", \"23-19=\"),\n    @(\"67-15=\", \"17+14=\"),\n    @(\"71-44=\", \"0+70=\"),\n    @(\"82-75=\", \"69-56=\"),\n    @(\"2+10=\", \"1+90=\"),\n    @(\"51+41=\", \"83-44=\"),\n    @(\"16+45=\", \"90-53=\"),\n    @(\"57-13=\", \"94-16=\"),\n    @(\"56-41=\", \"93-49=\"),\n    @(\"59+35=\", \"4+62=\"),\n    @(\"74-28=\", \"6+19=\"),\n    @(\"66-59=\", \"6+80=\"),\n    @(\"13+26=\", \"20+39=\"),\n    @(\"94-12=\", \"2+86=\"),\n    @(\"78-23=\", \"78-66=\"),\n    @(\"41-18=\", \"73-19=\"),\n    @(\"72+8=\", \"40+36=\"),\n    @(\"78+3=\", \"81-63=\"),\n    @(\"68-9=\", \"12+52=\"),\n    @(\"52+0=\", \"40-39=\"),\n)\n\n$wdFindContinue = 1\n$wdReplaceAll = 2\n\nforeach ($pair in $pairs) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Text = $pair[0]\n    $find.Replacement.ClearFormatting()\n    $find.Replacement.Text = $pair[1]\n    $find.Execute($pair[0], $false, $false, $false, $false, $false, $true, $wdFindContinue, $false, $pair[1], $wdReplaceAll)\n}\n"}
